# Added handler for Patient Status
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PatientStatus")

# Replace the existing status list with the new Patient Status handler codes.
$ws.Range("A1:B7").Value = $null

$ws.Range("A1").Value = "id"
$ws.Range("B1").Value = "status"

$ws.Range("A2").Value = 1
$ws.Range("B2").Value = "1 Preop"

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "2 For Exam"

$ws.Range("A5").Value = 4
$ws.Range("B5").Value = "4 Following"

$ws.Range("A6").Value = 5
$ws.Range("B6").Value = "5 Post-Op"

$ws.Range("A4").Value = 3
$ws.Range("B4").Value = "3 No Contact"

$ws.Range("A7").Value = 6
$ws.Range("B7").Value = "Deceased Preop"

$ws.Activate() | Out-Null
$ws.Range("B7").Select() | Out-Null
